# Append a new row (2025-04-19) to each of the 9 price-history sheets,
# carrying the previous day's price forward (sheet-specific new price for
# USD_CNY, which already happens to match 7.3173).
#
# Each sheet's data lives in columns A (Date, text) / B (Price, text) and
# the last populated row is 48. We need to add row 49 with:
#   A49 = "2025-04-19"   (plain text, not an Excel date serial)
#   B49 = <new price>    (plain text, not a number)
#
# Plain `Range.Value = "2025-04-19"` / `"39.5"` would be auto-coerced by
# Excel into a date serial / numeric value (and forcing text via
# NumberFormat="@" leaves a new cell style behind that isn't in the
# target). Instead we:
#   1. Copy row 48 down into row 49 (preserves the existing text cell
#      type/format with no style change).
#   2. Build the literal date/price strings with a throwaway formula
#      (e.g. ="2025-04-19") in a scratch cell far out of the used range,
#      copy that computed (text) result, and paste-special *values only*
#      over the cells that need to change. This yields a plain text value
#      with no extra number formatting / style footprint.
#   3. Clear the scratch cell again so it leaves no trace.

$wb = $excel.ActiveWorkbook

$newDate = "2025-04-19"

$sheets = @(
    @{ Name = "N-Dense";                   Price = "39.5" },
    @{ Name = "N-Type";                    Price = "40" },
    @{ Name = "N-type Wafer";              Price = "1.23" },
    @{ Name = "Cell Topcon 183mm";         Price = "0.293" },
    @{ Name = "Module Topcon 183mm";       Price = "0.09" },
    @{ Name = "Silver Rear_side";          Price = "5,329" },
    @{ Name = "Silver Busbar front-side";  Price = "7,977" },
    @{ Name = "Silver finger front-side";  Price = "8,027" },
    @{ Name = "USD_CNY";                   Price = "7.3173" }
)

foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.Name)

    $srcRow = 48
    $dstRow = 49

    # 1) Copy the last row down as-is (keeps text cell type, no new style).
    $ws.Range("A" + $srcRow + ":B" + $srcRow).Copy($ws.Range("A" + $dstRow + ":B" + $dstRow))

    # 2) Write the new date as literal text via a scratch formula cell.
    $scratch = $ws.Range("Z1000")
    $scratch.Formula = "=""" + $newDate + """"
    $scratch.Copy()
    $ws.Range("A" + $dstRow).PasteSpecial(-4163)

    # 3) Write the new price as literal text the same way.
    $scratch.Formula = "=""" + $s.Price + """"
    $scratch.Copy()
    $ws.Range("B" + $dstRow).PasteSpecial(-4163)

    # 4) Clean up the scratch cell so it leaves no trace.
    $scratch.Clear()
}
